$p = $ppt.ActivePresentation
$d = $p.Designs.Item(1)
$master = $d.SlideMaster
Write-Output "Master Name: $($master.Name)"
$theme = $master.Theme
Write-Output "Theme: $theme"
if ($theme -ne $null) {
  Write-Output "Theme.Name: $($theme.Name)"
}
